$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B59").Value = 44280
$ws.Range("B59").NumberFormat = "m/d/yyyy"

$ws.Range("C59").Value = 0.33333333333333331
$ws.Range("D59").Value = 0.57291666666666663

$ws.Range("F59").Value = "I-431"
$ws.Range("G59").Value = "Exercice"
$ws.Range("H59").Value = "Rapport"
$ws.Range("I59").Value = "CPNV"
$ws.Range("J59").Value = "M.Viret nous a donner un Rapport déjà fait en partie et nous avons du proposer quoi enlever / ajouter"
$ws.Range("K59").Value = "Oui"

$ws.Range("B60").Value = 44280
$ws.Range("C60").Value = 0.375
$ws.Range("D60").Value = 0.4375

$ws.Range("F60").Value = "Ma-20"
$ws.Range("G60").Value = "Code"
$ws.Range("H60").Value = "Score"
$ws.Range("I60").Value = "CPNV"
$ws.Range("J60").Value = "J'ai chercher comment enregistrer une variable dans un fichier texte a part et j'ai crée une fonction que fait ca"
$ws.Range("K60").Value = "oui"

$ws.Range("L58").Select()
